$wb = $excel.ActiveWorkbook

$wsLog = $wb.Worksheets.Item("Bug Log")
$wsMetrics = $wb.Worksheets.Item("Bug Metrics")

# Fill in row 8 of the Bug Log sheet with the new bug report
$wsLog.Range("A8").Value = 7
$wsLog.Range("B8").Value = 3
$wsLog.Range("C8").Value = "a1171e265566f2979ce323ba65d73b4ea24b24a2"
$wsLog.Range("D8").Value = "Reset Password"
$wsLog.Range("E8").Value = "Unable to retrieve the correct user from the database"
$wsLog.Range("F8").Value = "Error in the retrieval method for tutors"
$wsLog.Range("G8").Value = (Get-Date -Year 2018 -Month 8 -Day 29)
$wsLog.Range("H8").Value = "Zang Yu"
$wsLog.Range("I8").Value = (Get-Date -Year 2018 -Month 8 -Day 2)
$wsLog.Range("J8").Value = 0.5
$wsLog.Range("K8").Value = "High"
$wsLog.Range("L8").Value = 5
$wsLog.Range("M8").Value = 5
$wsLog.Range("N8").Value = "Resolved"
$wsLog.Range("O8").Value = "Corrected the retrieval method by modifying the query"

$wsLog.Select()
$wsLog.Range("O9").Select()

$wb.Windows.Item(1).ActiveSheet
